# Refresh the crypto price/volume snapshot (GitHub Actions daily update).
# Only the cells listed below changed between runs; every other cell (Coin
# name/index/link that didn't move) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '61.963.79'
$ws.Range('E2').Value = '  +4.77%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.084.17'
$ws.Range('E3').Value = '  +3.37%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.02%  '

# Row 5: BNB
$ws.Range('D5').Value = '''580.52'
$ws.Range('E5').Value = '  +3.18%  '

# Row 6: Solana
$ws.Range('D6').Value = '''142.16'
$ws.Range('E6').Value = '  +2.73%  '

# Row 7: USDC
$ws.Range('E7').Value = '  -0.18%  '

# Row 8: LidoStakedEther
$ws.Range('D8').Value = '3.072.78'
$ws.Range('E8').Value = '  +3.38%  '

# Row 9: XRP
$ws.Range('D9').Value = '''0.527'
$ws.Range('E9').Value = '  +1.25%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.140'
$ws.Range('E10').Value = '  +5.78%  '

# Row 11: Toncoin
$ws.Range('D11').Value = '''5.70'
$ws.Range('E11').Value = '  +11.16%  '

# Row 12: Cardano
$ws.Range('D12').Value = '''0.466'
$ws.Range('E12').Value = '  +2.79%  '

# Row 13: ShibaInu
$ws.Range('D13').Value = '''0.0000241'
$ws.Range('E13').Value = '  +4.94%  '

# Row 14: Avalanche
$ws.Range('D14').Value = '''35.34'
$ws.Range('E14').Value = '  +4.88%  '

# Row 15: TRON
$ws.Range('D15').Value = '''0.123'
$ws.Range('E15').Value = '  +0.18%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range('D16').Value = '3.590.51'
$ws.Range('E16').Value = '  +3.22%  '

# Row 17: Polkadot
$ws.Range('D17').Value = '''7.25'
$ws.Range('E17').Value = '  +1.71%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '3.079.20'
$ws.Range('E18').Value = '  +3.02%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '61.874.20'
$ws.Range('E19').Value = '  +4.64%  '

# Row 20: BitcoinCash
$ws.Range('D20').Value = '''447.66'
$ws.Range('E20').Value = '  +4.58%  '

# Row 21: Chainlink
$ws.Range('D21').Value = '''13.93'
$ws.Range('E21').Value = '  +2.51%  '

# Row 22: Polygon
$ws.Range('D22').Value = '''0.732'
$ws.Range('E22').Value = '  +2.47%  '

# Row 23: Uniswap
$ws.Range('E23').Value = '  +4.89%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range('D24').Value = '''13.79'
$ws.Range('E24').Value = '  +2.83%  '

# Row 25: Litecoin
$ws.Range('D25').Value = '''81.95'
$ws.Range('E25').Value = '  +1.56%  '

# Row 26: Dai
$ws.Range('E26').Value = '  +0.26%  '

# Row 27: ImmutableX
$ws.Range('E27').Value = '  +5.27%  '

# Row 28: FirstDigitalUSD
$ws.Range('E28').Value = '  -0.32%  '

# Row 29: PancakeSwap
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '''8.28'
$ws.Range('E29').Value = '  +7.48%  '

# Row 30: RenderToken
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''2.67'
$ws.Range('E30').Value = '  +5.11%  '

# Row 31: NEARProtocol
$ws.Range('D31').Value = '''6.80'
$ws.Range('E31').Value = '  +11.82%  '

# Row 32: Hedera
$ws.Range('D32').Value = '''0.111'
$ws.Range('E32').Value = '  +12.78%  '

# Row 33: EthereumClassic
$ws.Range('D33').Value = '''26.85'
$ws.Range('E33').Value = '  +4.50%  '

# Row 34: Mantle
$ws.Range('E34').Value = '  +4.27%  '

# Row 35: PEPE
$ws.Range('D35').Value = '0.0₃0794'
$ws.Range('E35').Value = '  +2.91%  '

# Row 36: Filecoin
$ws.Range('E36').Value = '  +3.61%  '

# Row 37: Stacks
$ws.Range('E37').Value = '  +5.40%  '

# Row 38: OKB
$ws.Range('D38').Value = '''50.21'
$ws.Range('E38').Value = '  +1.92%  '

# Row 39: dogwifhat
$ws.Range('E39').Value = '  +9.61%  '

# Row 40: Cosmos
$ws.Range('D40').Value = '''8.80'
$ws.Range('E40').Value = '  +2.07%  '

# Row 41: Bittensor
$ws.Range('D41').Value = '''422.09'
$ws.Range('E41').Value = '  +5.43%  '

# Row 42: VeChain
$ws.Range('D42').Value = '''0.0371'
$ws.Range('E42').Value = '  +5.60%  '

# Row 43: Maker
$ws.Range('D43').Value = '2.896.45'
$ws.Range('E43').Value = '  +4.50%  '

# Row 44: TheGraph
$ws.Range('D44').Value = '''0.274'
$ws.Range('E44').Value = '  +9.29%  '

# Row 45: Kaspa
$ws.Range('E45').Value = '  +0.78%  '

# Row 46: Fetch.AI
$ws.Range('D46').Value = '''2.14'
$ws.Range('E46').Value = '  +7.25%  '

# Row 47: USDe
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '''35.18'
$ws.Range('E47').Value = '  +3.99%  '

# Row 48: Arweave
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '''0.999'
$ws.Range('E48').Value = '  -0.01%  '

# Row 49: Monero
$ws.Range('D49').Value = '''123.79'
$ws.Range('E49').Value = '  +2.12%  '

# Row 50: Stellar
$ws.Range('E50').Value = '  +1.15%  '

# Row 51: InjectiveProtocol
$ws.Range('D51').Value = '''24.18'
$ws.Range('E51').Value = '  +3.32%  '
